$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right before the current row 635, shifting the
# existing rows 635:676 down to 637:678.
$ws.Rows("635:636").Insert()

# New row 635: Femacal de La Calera - Piña Caramelo - Primera
$ws.Range("A635").Value = 3
$ws.Range("B635").Value = "Femacal de La Calera"
$ws.Range("C635").Value = "Coquimbo"
$ws.Range("D635").Value = 44746
$ws.Range("E635").Value = 5
$ws.Range("F635").Value = "Fruta"
$ws.Range("G635").Value = 100108
$ws.Range("H635").Value = "Tropicales y subtropicales"
$ws.Range("I635").Value = 100108005
$ws.Range("J635").Value = "Piña"
$ws.Range("K635").Value = "Caramelo"
$ws.Range("L635").Value = "Primera"
$ws.Range("M635").Value = 108
$ws.Range("N635").Value = 22000
$ws.Range("O635").Value = 22000
$ws.Range("P635").Value = 22000
$ws.Range("Q635").Value = "$/caja 12 unidades"
$ws.Range("R635").Value = "Ecuador"
$ws.Range("S635").Value = 1833
$ws.Range("T635").Value = 12

# New row 636: Femacal de La Calera - Piña Caramelo - Segunda
$ws.Range("A636").Value = 3
$ws.Range("B636").Value = "Femacal de La Calera"
$ws.Range("C636").Value = "Coquimbo"
$ws.Range("D636").Value = 44746
$ws.Range("E636").Value = 5
$ws.Range("F636").Value = "Fruta"
$ws.Range("G636").Value = 100108
$ws.Range("H636").Value = "Tropicales y subtropicales"
$ws.Range("I636").Value = 100108005
$ws.Range("J636").Value = "Piña"
$ws.Range("K636").Value = "Caramelo"
$ws.Range("L636").Value = "Segunda"
$ws.Range("M636").Value = 54
$ws.Range("N636").Value = 22000
$ws.Range("O636").Value = 22000
$ws.Range("P636").Value = 22000
$ws.Range("Q636").Value = "$/caja 14 unidades"
$ws.Range("R636").Value = "Ecuador"
$ws.Range("S636").Value = 1571
$ws.Range("T636").Value = 14
